# conditionAlwayFalse-template.docx edit:
#  - Move the "_GoBack" bookmark from paragraph 2 (inside the IF field code)
#    to the very start of paragraph 1 (where the cursor/selection was left
#    when the document was last saved).
#  - Word's background spell-checker re-ran over the three sentences that
#    contain the word "demonstration"/"paragraph" and wrapped each of those
#    words with <w:proofErr w:type="spellStart"/>...<w:proofErr w:type="spellEnd"/>,
#    which also splits the surrounding <w:r> runs at the word boundaries.
#
# InsertXML() replaces the exact contents of the Range it is called on, so
# each affected paragraph's Range is rewritten wholesale with its new,
# finer-grained run/bookmark/proofErr structure (this also conveniently
# sidesteps this host's Bookmarks.Add() edge case for a zero-length range
# that sits at absolute document position 0).

$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)
$p1.Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`" w:rsidR=`"00BA34AE`" w:rsidRDefault=`"00BA34AE`" w:rsidP=`"00F5495F`"><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/><w:r><w:t xml:space=`"preserve`">Basic </w:t></w:r><w:r w:rsidR=`"002A1F2A`"><w:t>if</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>demonstration</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t> :</w:t></w:r></w:p>") | Out-Null

$p2 = $d.Paragraphs(2)
$p2.Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`" w:rsidR=`"002A1F2A`" w:rsidRDefault=`"002A1F2A`" w:rsidP=`"002A1F2A`"><w:pPr><w:tabs><w:tab w:val=`"left`" w:pos=`"3119`"/></w:tabs></w:pPr><w:r><w:fldChar w:fldCharType=`"begin`"/></w:r><w:r><w:instrText xml:space=`"preserve`">m:if </w:instrText></w:r><w:r w:rsidR=`"00F92982`"><w:instrText>self.oclIsKindOf(ecore::EC</w:instrText></w:r><w:r w:rsidR=`"00141F2F`"><w:instrText>l</w:instrText></w:r><w:r w:rsidR=`"00F92982`"><w:instrText>ass)</w:instrText></w:r><w:r><w:fldChar w:fldCharType=`"end`"/></w:r></w:p>") | Out-Null

$p3 = $d.Paragraphs(3)
$p3.Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`" w:rsidR=`"002A1F2A`" w:rsidRDefault=`"000321D3`" w:rsidP=`"002A1F2A`"><w:pPr><w:tabs><w:tab w:val=`"left`" w:pos=`"3119`"/></w:tabs></w:pPr><w:r><w:t xml:space=`"preserve`">The </w:t></w:r><w:r w:rsidR=`"006E6505`"><w:t>THEN</w:t></w:r><w:r w:rsidR=`"002A1F2A`"><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>paragraph</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t>.</w:t></w:r></w:p>") | Out-Null

$p5 = $d.Paragraphs(5)
$p5.Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`" w:rsidR=`"00C52979`" w:rsidRDefault=`"00BA34AE`" w:rsidP=`"00F5495F`"><w:r><w:t>En</w:t></w:r><w:r w:rsidR=`"006B5B12`"><w:t>d</w:t></w:r><w:r><w:t xml:space=`"preserve`"> of </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>demonstration</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t>.</w:t></w:r></w:p>") | Out-Null
